$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("691_3")

# Row 3 corresponds to "Петров". Fix bug for delete rebuke:
# move the rebuke mark from column C (НК) to column E (старшина)
$ws.Range("C3").Value = 0
$ws.Range("D3").Value = 0
$ws.Range("E3").Value = 1
